# Generate Report for Handoff
#
# The handoff xliff files were (re-)generated, so the status moves from
# "In Translation" to "Ready for handoff" and the associated generation /
# handoff timestamps are refreshed. Excel widens the "status" columns to
# fit the longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# Columns: A File Name | B Path And Name | C Extension | D Publish URL |
#          E zh-cn | F de-de | G Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-28 02:58:31"

# --- zh-cn sheet -------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-28 02:58:27"

# --- de-de sheet -------------------------------------------------------
# Column C = Status; its Latest Handoff Datetime (H2) mirrors the
# Overview "Latest HO Xliff Generate Date" value set above.
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-28 02:58:31"

# --- Column widths ------------------------------------------------------
# The status columns grow wider to accommodate "Ready for handoff"
# (target character width ~17.216).
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
